$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2: reorder "Recorded By" list, moving "System" to the end
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"

# Row 3: add an additional recorder and update the attendance count
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("H3").Value = "62/251"

# Average Attendance % figures updated to reflect the new count.
# A leading apostrophe forces these to be stored as literal text
# (matching the original cells) instead of Excel auto-converting
# "26.2%" into a numeric percentage value.
$ws.Range("L10").Value = "'26.2%"
$ws.Range("S15").Value = "'26.2%"

# The apostrophe entry flags the cells with a quote-prefix style; restore
# the original plain formatting by re-pasting the formats from a
# same-styled neighbour cell that already holds a literal percent string.
$ws.Range("L9").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 28: reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
